# Scheduled-runner style refresh of market-price columns (currentAveragePrice,
# currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ) across the ALC,
# ARM, BSM, CRP, CUL, GSM, LTW and WVR sheets. Values below mirror the
# upstream data-pull; a handful of rows had their profit cell become blank
# (cleared, not zeroed) where the source feed no longer returns a price.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 885.75
$ws.Range("I12").Value = 343
$ws.Range("K12").Value = 343
$ws.Range("M12").Value = -173

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 540.6667
$ws.Range("I33").Value = 488.46155
$ws.Range("J33").Value = 676.4
$ws.Range("K33").Value = 488.46155
$ws.Range("L33").Value = 676.4
$ws.Range("M33").Value = -259.46155
$ws.Range("N33").Value = -1134.4

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 286.66666
$ws.Range("I41").Value = 289
$ws.Range("J41").Value = 275
$ws.Range("K41").Value = 289
$ws.Range("L41").Value = 275
$ws.Range("M41").Value = 151
$ws.Range("N41").Value = -1155

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4850
$ws.Range("I51").Value = 4850
$ws.Range("K51").Value = 4850
$ws.Range("M51").Value = -4366

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H96").Value = 483.16666
$ws.Range("I96").Value = 483.16666
$ws.Range("K96").Value = 1449.49998
$ws.Range("M96").Value = -76.49998000000005

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 4989.2856
$ws.Range("J97").Value = 4989.2856
$ws.Range("L97").Value = 14967.8568
$ws.Range("N97").Value = -15959.8568

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 4262.3823
$ws.Range("I98").Value = 4240.788
$ws.Range("J98").Value = 4975
$ws.Range("K98").Value = 4240.788
$ws.Range("L98").Value = 4975
$ws.Range("M98").Value = -2742.788
$ws.Range("N98").Value = -7971

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 4262.3823
$ws.Range("I122").Value = 4240.788
$ws.Range("J122").Value = 4975
$ws.Range("K122").Value = 12722.364
$ws.Range("L122").Value = 14925
$ws.Range("M122").Value = -10272.364
$ws.Range("N122").Value = -19825

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 3000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 9000
$ws.Range("N137").Value = -14100
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 5466.421
$ws.Range("I141").Value = 4619.5625
$ws.Range("K141").Value = 13858.6875
$ws.Range("M141").Value = -8678.6875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3737.7585
$ws.Range("I2").Value = 3284.3
$ws.Range("K2").Value = 3284.3
$ws.Range("M2").Value = -3171.3

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2337.9656
$ws.Range("I32").Value = 2064.3215
$ws.Range("K32").Value = 2064.3215
$ws.Range("M32").Value = -1777.3215

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4908645
$ws.Range("I61").Value = 7580815.5
$ws.Range("K61").Value = 7580815.5
$ws.Range("M61").Value = -7580603.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1082.3334
$ws.Range("I97").Value = 998
$ws.Range("K97").Value = 998
$ws.Range("M97").Value = -502

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 3737.7585
$ws.Range("I116").Value = 3284.3
$ws.Range("K116").Value = 3284.3
$ws.Range("M116").Value = -990.3000000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 3988.1538
$ws.Range("I122").Value = 3988.1538
$ws.Range("K122").Value = 11964.4614
$ws.Range("M122").Value = -9514.4614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 5864.2812
$ws.Range("I132").Value = 4876.696
$ws.Range("K132").Value = 14630.088
$ws.Range("M132").Value = -12100.088

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4908645
$ws.Range("I136").Value = 7580815.5
$ws.Range("K136").Value = 22742446.5
$ws.Range("M136").Value = -22739896.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3737.7585
$ws.Range("I3").Value = 3284.3
$ws.Range("K3").Value = 3284.3
$ws.Range("M3").Value = -3170.3

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 11766382
$ws.Range("I86").Value = 1776.2307
$ws.Range("J86").Value = 50001350
$ws.Range("K86").Value = 1776.2307
$ws.Range("L86").Value = 50001350
$ws.Range("M86").Value = -653.2307000000001
$ws.Range("N86").Value = -50003596

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 11766382
$ws.Range("I89").Value = 1776.2307
$ws.Range("J89").Value = 50001350
$ws.Range("K89").Value = 8881.1535
$ws.Range("L89").Value = 250006750
$ws.Range("M89").Value = -3265.1535
$ws.Range("N89").Value = -250017982

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 8249.833000000001
$ws.Range("J105").Value = 8999.666999999999
$ws.Range("L105").Value = 8999.666999999999
$ws.Range("N105").Value = -12493.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 4113.8667
$ws.Range("I107").Value = 3208.3076
$ws.Range("K107").Value = 3208.3076
$ws.Range("M107").Value = -1288.3076

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3499.8333
$ws.Range("J16").Value = 4068.1
$ws.Range("L16").Value = 4068.1
$ws.Range("N16").Value = -4642.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4368.722
$ws.Range("I31").Value = 3132.3215
$ws.Range("K31").Value = 3132.3215
$ws.Range("M31").Value = -2837.3215

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 0
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4368.722
$ws.Range("I34").Value = 3132.3215
$ws.Range("K34").Value = 3132.3215
$ws.Range("M34").Value = -2930.3215

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3201.25
$ws.Range("I62").Value = 2933
$ws.Range("J62").Value = 4006
$ws.Range("K62").Value = 2933
$ws.Range("L62").Value = 4006
$ws.Range("M62").Value = -2309
$ws.Range("N62").Value = -5254

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 3201.25
$ws.Range("I65").Value = 2933
$ws.Range("J65").Value = 4006
$ws.Range("K65").Value = 14665
$ws.Range("L65").Value = 20030
$ws.Range("M65").Value = -11545
$ws.Range("N65").Value = -26270

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2256.6667
$ws.Range("I105").Value = 1708
$ws.Range("J105").Value = 5000
$ws.Range("K105").Value = 1708
$ws.Range("L105").Value = 5000
$ws.Range("M105").Value = 39
$ws.Range("N105").Value = -8494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2931.111
$ws.Range("J107").Value = 4300
$ws.Range("L107").Value = 4300
$ws.Range("N107").Value = -8140

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 3499.8333
$ws.Range("J113").Value = 4068.1
$ws.Range("L113").Value = 4068.1
$ws.Range("N113").Value = -8408.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2287.2917
$ws.Range("I132").Value = 1432.5625
$ws.Range("K132").Value = 4297.6875
$ws.Range("M132").Value = -1767.6875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 3860.8076
$ws.Range("I134").Value = 2042.8125
$ws.Range("J134").Value = 6769.6
$ws.Range("K134").Value = 6128.4375
$ws.Range("L134").Value = 20308.8
$ws.Range("M134").Value = -3593.4375
$ws.Range("N134").Value = -25378.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("M16").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 1482.8334
$ws.Range("I13").Value = 174.25
$ws.Range("K13").Value = 174.25
$ws.Range("M13").Value = -35.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 828.5714
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H82").Value = 69328
$ws.Range("J82").Value = 69328
$ws.Range("L82").Value = 69328
$ws.Range("N82").Value = -70094

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H85").Value = 69328
$ws.Range("J85").Value = 69328
$ws.Range("L85").Value = 69328
$ws.Range("N85").Value = -71980

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2721.2415
$ws.Range("I102").Value = 2378
$ws.Range("J102").Value = 4368.8
$ws.Range("K102").Value = 2378
$ws.Range("L102").Value = 4368.8
$ws.Range("M102").Value = -756
$ws.Range("N102").Value = -7612.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2337.375
$ws.Range("I122").Value = 2337.375
$ws.Range("K122").Value = 7012.125
$ws.Range("M122").Value = -4562.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 6090.5415
$ws.Range("I132").Value = 5627.619
$ws.Range("J132").Value = 9331
$ws.Range("K132").Value = 16882.857
$ws.Range("L132").Value = 27993
$ws.Range("M132").Value = -14352.857
$ws.Range("N132").Value = -33053

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 17695.076
$ws.Range("I93").Value = 1003.9
$ws.Range("K93").Value = 1003.9
$ws.Range("M93").Value = 244.1

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 9093.579
$ws.Range("I132").Value = 10423
$ws.Range("J132").Value = 6213.1665
$ws.Range("K132").Value = 31269
$ws.Range("L132").Value = 18639.4995
$ws.Range("M132").Value = -28739
$ws.Range("N132").Value = -23699.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 10002699
$ws.Range("I2").Value = 20002398
$ws.Range("K2").Value = 20002398
$ws.Range("M2").Value = -20002286

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9998.5
$ws.Range("J62").Value = 9998.5
$ws.Range("L62").Value = 9998.5
$ws.Range("N62").Value = -11246.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 9998.5
$ws.Range("J65").Value = 9998.5
$ws.Range("L65").Value = 49992.5
$ws.Range("N65").Value = -56232.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4986.936
$ws.Range("I132").Value = 4492.0625
$ws.Range("K132").Value = 13476.1875
$ws.Range("M132").Value = -10946.1875
